# Preparatory work for PrOD combination tables: add ombitasvir, paritaprevir,
# dasabuvir rows, pulled in with the same formatting as the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert three new rows (11-13), copying row 10's formatting so the
#     abbreviation column keeps the same style (s="3") used by rows 8-10. ---
$ws.Range("A10:E10").Copy()
$ws.Range("A11:E11").Insert(-4121)
$ws.Range("A10:E10").Copy()
$ws.Range("A12:E12").Insert(-4121)
$ws.Range("A10:E10").Copy()
$ws.Range("A13:E13").Insert(-4121)

# The insert operations pushed the old trailing formatted row (16) down to
# row 19 - move it back up to row 16 and tidy up the now-empty row 19.
$ws.Range("A19:B19").Cut($ws.Range("A16:B16"))
$ws.Range("A19:B19").Clear()

# --- Fill in the new drug data. Shared strings are created in this exact
#     write order: A11, A12, A13, E11, B11, B12, B13, E12, E13. ---
$ws.Range("A11").Value = "ombitasvir"
$ws.Range("A12").Value = "paritaprevir"
$ws.Range("A13").Value = "dasabuvir"

$ws.Range("E11").Value = "ABT-267"

$ws.Range("B11").Value = "OBV"
$ws.Range("B12").Value = "PTV"
$ws.Range("B13").Value = "DSV"

$ws.Range("E12").Value = "ABT-450"
$ws.Range("E13").Value = "ABT-333"

# category / producer columns reuse strings already present in the sheet.
$ws.Range("C11").Value = "NS5A" + [char]0x00A0 + "inhibitors"
$ws.Range("C12").Value = "NS3/4A" + [char]0x00A0 + "protease inhibitors"
$ws.Range("C13").Value = "NS5B RNA polymerase inhibitors"

$ws.Range("D11").Value = "Abbvie"
$ws.Range("D12").Value = "Abbvie"
$ws.Range("D13").Value = "Abbvie"

# --- Selection left on E15 after the edit ---
$ws.Range("E15").Select() | Out-Null
